$wb = $excel.ActiveWorkbook

# --- new_vars sheet: add two new rows of data ---
$ws = $wb.Worksheets.Item("new_vars")

# Row 3: language / expressive vocabulary / WOLD Expressive Vocabulary
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "language"
$ws.Range("C3").Value = "expressive vocabulary"
$ws.Range("D3").Value = "WOLD Expressive Vocabulary"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = "wold_vcb_raw_f8"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = "WP1_DataExtraction.Rmd"

# Row 4: executive function / attention / TEACh Attentional Control
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "executive function"
$ws.Range("C4").Value = "attention"
$ws.Range("D4").Value = "TEACh Attentional Control"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = "teach_ctr_diff_f8"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = "WP1_DataExtraction.Rmd"

# Column widths for A:D on new_vars
$ws.Range("A:A").ColumnWidth = 12.28515625
$ws.Range("B:B").ColumnWidth = 14.85546875
$ws.Range("C:C").ColumnWidth = 17.140625
$ws.Range("D:D").ColumnWidth = 28.42578125

# Selection on new_vars ends at F4
$ws.Range("F4").Select()

# --- metadata sheet: reset selection to A1 ---
$meta = $wb.Worksheets.Item("metadata")
$meta.Range("A1").Select()

# keep new_vars the active sheet/tab
$ws.Activate()
